$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (non price) column updates: Coin name / Link / Volume% columns ---
# These values are not numeric-looking, so a direct .Value assignment is safe
# and will not be auto-coerced to a number by Excel.
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("E16").Value = "  +1.90%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("E36").Value = "  +4.55%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("E42").Value = "  +4.34%  "
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("E51").Value = "  +3.45%  "

# --- Price (D) column updates ---
# These values can look numeric (e.g. "311.17"), and Excel COM will silently
# convert a plain .Value assignment of such a string into a Number, which would
# both change the stored cell type (Text -> Number) and risk float round-off.
# To keep them as literal text (matching the original inlineStr cells) without
# perturbing the cell style, stage each one as a text-producing formula, then
# convert the whole block to static values via Copy + PasteSpecial(xlPasteValues).
$ws.Range("D2").Formula = "=""27.421.32"""
$ws.Range("D3").Formula = "=""1.862.88"""
$ws.Range("D5").Formula = "=""311.17"""
$ws.Range("D7").Formula = "=""0.4773"""
$ws.Range("D8").Formula = "=""0.3759"""
$ws.Range("D9").Formula = "=""0.07314"""
$ws.Range("D10").Formula = "=""0.9335"""
$ws.Range("D11").Formula = "=""20.62"""
$ws.Range("D12").Formula = "=""0.07807"""
$ws.Range("D13").Formula = "=""1.915.00"""
$ws.Range("D14").Formula = "=""5.425"""
$ws.Range("D15").Formula = "=""6.553"""
$ws.Range("D16").Formula = "=""90.48"""
$ws.Range("D18").Formula = "=""0.000008877"""
$ws.Range("D19").Formula = "=""1.010"""
$ws.Range("D20").Formula = "=""27.458.65"""
$ws.Range("D21").Formula = "=""14.67"""
$ws.Range("D22").Formula = "=""5.111"""
$ws.Range("D23").Formula = "=""10.68"""
$ws.Range("D24").Formula = "=""1.938"""
$ws.Range("D25").Formula = "=""155.45"""
$ws.Range("D26").Formula = "=""18.45"""
$ws.Range("D27").Formula = "=""2.013"""
$ws.Range("D28").Formula = "=""115.28"""
$ws.Range("D29").Formula = "=""4.941"""
$ws.Range("D30").Formula = "=""0.08887"""
$ws.Range("D31").Formula = "=""3.323"""
$ws.Range("D32").Formula = "=""1.212"""
$ws.Range("D33").Formula = "=""4.592"""
$ws.Range("D34").Formula = "=""0.7515"""
$ws.Range("D35").Formula = "=""2.724"""
$ws.Range("D36").Formula = "=""0.02043"""
$ws.Range("D37").Formula = "=""1.114"""
$ws.Range("D38").Formula = "=""0.05256"""
$ws.Range("D39").Formula = "=""2.985"""
$ws.Range("D40").Formula = "=""0.5318"""
$ws.Range("D41").Formula = "=""7.062"""
$ws.Range("D42").Formula = "=""8.570"""
$ws.Range("D45").Formula = "=""0.4797"""
$ws.Range("D47").Formula = "=""1.654"""
$ws.Range("D48").Formula = "=""102.75"""
$ws.Range("D49").Formula = "=""67.16"""
$ws.Range("D50").Formula = "=""0.06077"""
$ws.Range("D51").Formula = "=""0.9187"""

$priceRange = $ws.Range("D2:D51")
$priceRange.Copy()
$priceRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false
